# Apply the LinuxForHealth re-brand / version-bump edit described by the
# commit "Deploying to gh-pages from @ LinuxForHealth/alvearie-fhir-ig@...".
#
# Workbook layout:
#   Sheet 1 "Metadata" - simple Property/Value table (A:B)
#   Sheet 2 "Elements" - StructureDefinition element table (A:AJ)

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$elements = $wb.Worksheets.Item("Elements")

# --- Metadata sheet -------------------------------------------------
# URL: ibm.com -> linuxforhealth.org
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/measure-report-evidence-value"

# Version bump 7.0.0 -> 8.0.0
$meta.Range("B3").Value = "8.0.0"

# Date refreshed
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"

# Publisher renamed
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet ---------------------------------------------------
# The "Fixed Value" for Extension.url mirrors the canonical extension URL
# above, so it moves to linuxforhealth.org too.
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/measure-report-evidence-value"

# The root Extension row no longer carries the ele-1/ext-1 constraint text
# in its "Constraint(s)" column (it stays correctly attached to the
# Extension.extension row only).
$elements.Range("AI2").Value = ""

Write-Host "edits applied"
